$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# New matchup rows for fall 22 week 10 (rows 1385-1408), columns A-D:
# Player_1, Points_1, Player_2, Points_2
$data = @(
    @(4,2,4,1),
    @(3,2,2,1),
    @(5,0,4,2),
    @(3,0,3,3),
    @(4,0,5,2),
    @(4,3,4,0),
    @(5,0,6,2),
    @(4,0,7,3),
    @(5,3,3,0),
    @(5,3,3,0),
    @(6,0,5,2),
    @(5,0,4,3),
    @(5,2,6,0),
    @(4,2,4,0),
    @(2,3,3,0),
    @(4,1,5,2),
    @(4,2,4,1),
    @(6,2,6,0),
    @(5,2,5,0),
    @(3,3,3,0),
    @(3,0,5,3),
    @(3,3,2,0),
    @(3,1,4,2),
    @(5,1,5,2)
)

$startRow = 1385
for ($i = 0; $i -lt $data.Count; $i++) {
    $r = $startRow + $i
    $row = $data[$i]
    for ($c = 0; $c -lt $row.Count; $c++) {
        $ws.Cells.Item($r, $c + 1).Value = $row[$c]
    }
}

$endRow = $startRow + $data.Count - 1
$nextRow = $endRow + 1

# Update scroll position / selection to reflect the new end of data,
# matching where Excel would land after entering this block of rows.
$excel.ActiveWindow.ScrollRow = $startRow
$selCell = "A" + $nextRow
$ws.Range($selCell).Select()

